$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.674.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -7.49%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.544.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -4.11%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.17%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'299.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -3.77%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'93.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -4.98%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -4.13%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.02%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -5.57%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'36.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -6.92%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0802"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -5.29%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'7.70"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -5.02%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +5.21%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.929.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -3.99%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.543.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -4.23%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.871"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -5.85%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'14.21"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -5.05%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'42.709.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -7.52%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'12.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.12%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.0₃0980"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -3.98%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -3.76%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'71.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -4.47%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'256.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -9.34%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -4.77%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("B25").Value = "'ImmutableX"
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = "'2.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -5.91%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("B26").Value = "'EthereumClassic"
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = "'29.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -3.62%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -0.11%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'10.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -5.26%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'36.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -4.41%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -5.63%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'5.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -5.23%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'152.46"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.62%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -8.15%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -2.54%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -10.26%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.0792"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -6.15%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -7.49%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.120"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -3.47%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'17.06"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +7.83%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'23.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.37%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -6.00%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -4.38%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -5.54%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'2.079.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -3.16%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -0.01%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'9.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.05%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +3.76%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'84.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -10.09%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'2.787.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -3.93%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -6.24%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -5.30%  "
$ws.Range("E51").Style = "Normal"
